$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11
$ws.Cells.Item(11, 1).Value = 9816.16
$ws.Cells.Item(11, 2).Value = 9822.0499999999993
$ws.Cells.Item(11, 3).Value = 283.47000000000003
$ws.Cells.Item(11, 4).Value = 283.29000000000002
$ws.Cells.Item(11, 5).Value = $false
$ws.Cells.Item(11, 6).Value = -0.06
$ws.Cells.Item(10, 7).Copy($ws.Cells.Item(11, 7))
$ws.Cells.Item(11, 7).Value = 42613.765590277777
$ws.Cells.Item(11, 8).Value = $false

# Row 12
$ws.Cells.Item(12, 1).Value = 9764.1299999999992
$ws.Cells.Item(12, 2).Value = 9816.16
$ws.Cells.Item(12, 3).Value = 282.39
$ws.Cells.Item(12, 4).Value = 280.89
$ws.Cells.Item(12, 5).Value = $false
$ws.Cells.Item(12, 6).Value = -0.53
$ws.Cells.Item(10, 7).Copy($ws.Cells.Item(12, 7))
$ws.Cells.Item(12, 7).Value = 42614.672997685186
$ws.Cells.Item(12, 8).Value = $false

# Row 13
$ws.Cells.Item(13, 1).Value = 9792.4500000000007
$ws.Cells.Item(13, 2).Value = 9764.1299999999992
$ws.Cells.Item(13, 3).Value = 280.62
$ws.Cells.Item(13, 4).Value = 281.44
$ws.Cells.Item(13, 5).Value = $false
$ws.Cells.Item(13, 6).Value = 0.28999999999999998
$ws.Cells.Item(10, 7).Copy($ws.Cells.Item(13, 7))
$ws.Cells.Item(13, 7).Value = 42615.750196759262
$ws.Cells.Item(13, 8).Value = $true
